# RMs scheduler architectures, updated INP common patterns and aspects of
# interests to RMs.
#
# On the "NetChain: logical communication pattern" slide (sldId 265), a
# handful of shapes were brought in front of everything else: two curved
# connector arrows (ids 28 and 30), the green "1. Direct communication..."
# callout textbox (id 38) together with its feeding connector (id 50), and
# the green connector (id 36). PowerPoint records this as each shape's
# z-order changing ("ord") while nothing else about them (position, text,
# style, connections) is touched, so we reproduce it with ZOrder(ppBringToFront)
# applied to each shape in turn, in the same relative order they ended up in.

$p = $ppt.ActivePresentation

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Locate the "NetChain: logical communication pattern" slide (sldId 265)
# by its title text rather than a hard-coded index.
$slide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.HasTitle) {
        $titleText = $candidate.Shapes.Title.TextFrame.TextRange.Text
        if ($titleText -eq "NetChain: logical communication pattern") {
            $slide = $candidate
            break
        }
    }
}

$msoBringToFront = 0

$orderedIds = @(28, 30, 38, 50, 36)
foreach ($id in $orderedIds) {
    $shape = Get-ShapeById $slide $id
    $shape.ZOrder($msoBringToFront)
}
